$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial value that was updated from
# 45172 (2023-09-03) to 45175 (2023-09-06) for every data row (rows 2-484).
$ws.Range("C2:C484").Value = 45175
